# Actualización automática 2025-10-23 16:30:08
#
# Applies the monthly sales-report refresh to the three worksheets:
#   - "VENTAS POR GRUPO"      (per-client sales by product group)
#   - "VENTA MENSUAL"         (per-client sales by month, incl. "octubre")
#   - "CUMPLIMIENTO MENSUAL"  (aggregated budget-compliance summary)
#
# All of the cells in this workbook hold literal cached values (there are
# no live formulas), so every downstream total/percentage/count that
# depends on the updated source cells must be written explicitly as well.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# -----------------------------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsVentasGrupo.Range("M4").Value  = 5125.26
$wsVentasGrupo.Range("L5").Value  = 2838.61
$wsVentasGrupo.Range("M5").Value  = 275.57
$wsVentasGrupo.Range("K16").Value = 319.68
$wsVentasGrupo.Range("L16").Value = 2244.24
$wsVentasGrupo.Range("M16").Value = 72.53
$wsVentasGrupo.Range("M24").Value = 6964.08
$wsVentasGrupo.Range("K36").Value = 3928.08
$wsVentasGrupo.Range("M37").Value = 3045.42

# Row 56 holds "<count> de 54" labels -- the count of clients (rows 2-55)
# with a strictly-positive value in that column. K/L/M each gained one
# more qualifying client because of the edits above.
$wsVentasGrupo.Range("K56").Value = "7 de 54"
$wsVentasGrupo.Range("L56").Value = "4 de 54"
$wsVentasGrupo.Range("M56").Value = "11 de 54"

# -----------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# -----------------------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsVentaMensual.Range("F4").Value  = 6745.25
$wsVentaMensual.Range("F5").Value  = 3114.18
$wsVentaMensual.Range("F16").Value = 4074.47
$wsVentaMensual.Range("F24").Value = 7750.68
$wsVentaMensual.Range("F36").Value = 13855.8
$wsVentaMensual.Range("F37").Value = 4896.77

# Row 60 is the column total for "octubre" (F2:F59).
$wsVentaMensual.Range("F60").Value = 63399.43

# -----------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# -----------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 10: PANELES DECORATIVOS  -- VENTA / POR CUMPLIR / CUMPLIMIENTO
$wsCumplimiento.Range("D10").Value = 8272.49
$wsCumplimiento.Range("E10").Value = -4391.410164656079
$wsCumplimiento.Range("F10").Value = 2.131491840148384

# Row 11: PIEDRA SINTERIZADA
$wsCumplimiento.Range("D11").Value = 8800.24
$wsCumplimiento.Range("E11").Value = 3030.76
$wsCumplimiento.Range("F11").Value = 0.7438289240131857

# Row 12: PORCELANATO
$wsCumplimiento.Range("D12").Value = 32542.85
$wsCumplimiento.Range("E12").Value = 20120.27
$wsCumplimiento.Range("F12").Value = 0.6179438286223832

# Row 14: TOTAL
$wsCumplimiento.Range("D14").Value = 61304.94
$wsCumplimiento.Range("E14").Value = 37711.56661190614
$wsCumplimiento.Range("F14").Value = 0.619138587066941
